$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Uzbekistan moves ahead of Bulgaria (rows 81/82), Uzbekistan stats updated ---
$ws.Cells.Item(81, 1).Value = "Uzbekistan"
$ws.Cells.Item(81, 2).Value = 669
$ws.Cells.Item(81, 3).Value = 45
$ws.Cells.Item(81, 4).Value = 42
$ws.Cells.Item(81, 5).Value = 624
$ws.Cells.Item(81, 6).Value = 8
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 3

$ws.Cells.Item(82, 1).Value = "Bulgaria"
$ws.Cells.Item(82, 2).Value = 635
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 54
$ws.Cells.Item(82, 5).Value = 556
$ws.Cells.Item(82, 6).Value = 33
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 25

# --- Sri Lanka moves ahead of Mayotte (rows 116/117), Sri Lanka stats updated ---
$ws.Cells.Item(116, 1).Value = "Sri Lanka"
$ws.Cells.Item(116, 2).Value = 197
$ws.Cells.Item(116, 3).Value = 7
$ws.Cells.Item(116, 4).Value = 54
$ws.Cells.Item(116, 5).Value = 136
$ws.Cells.Item(116, 6).Value = 5
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 7

$ws.Cells.Item(117, 1).Value = "Mayotte"
$ws.Cells.Item(117, 2).Value = 191
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = 50
$ws.Cells.Item(117, 5).Value = 139
$ws.Cells.Item(117, 6).Value = 4
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 2

# --- Guatemala moves ahead of Brunei/Paraguay/Gibraltar (rows 124-127) ---
$ws.Cells.Item(124, 1).Value = "Guatemala"
$ws.Cells.Item(124, 2).Value = 137
$ws.Cells.Item(124, 3).Value = 11
$ws.Cells.Item(124, 4).Value = 19
$ws.Cells.Item(124, 5).Value = 115
$ws.Cells.Item(124, 6).Value = 3
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 3

$ws.Cells.Item(125, 1).Value = "Brunei"
$ws.Cells.Item(125, 2).Value = 136
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 4).Value = 99
$ws.Cells.Item(125, 5).Value = 36
$ws.Cells.Item(125, 6).Value = 3
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 1

$ws.Cells.Item(126, 1).Value = "Paraguay"
$ws.Cells.Item(126, 2).Value = 133
$ws.Cells.Item(126, 3).Value = 4
$ws.Cells.Item(126, 4).Value = 18
$ws.Cells.Item(126, 5).Value = 109
$ws.Cells.Item(126, 6).Value = 1
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 6

$ws.Cells.Item(127, 1).Value = "Gibraltar"
$ws.Cells.Item(127, 2).Value = 127
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 69
$ws.Cells.Item(127, 5).Value = 58
$ws.Cells.Item(127, 6).Value = 1
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

# --- Jamaica stats updated (row 141, no reordering) ---
$ws.Cells.Item(141, 1).Value = "Jamaica"
$ws.Cells.Item(141, 2).Value = 65
$ws.Cells.Item(141, 3).Value = 2
$ws.Cells.Item(141, 4).Value = 13
$ws.Cells.Item(141, 5).Value = 48
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 4

# --- Update "last updated" timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 11 de Abril de 2020 a las 04:22"
